# Update the NIT-9010759661 "Estado de Cuenta" worksheet:
#  - Refresh the VALOR MORA total and the worker/period counters
#  - Replace the worker detail table with the new data set (part 1 of the
#    new estado de cuenta), which has 7 data rows instead of 9

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data table currently spans rows 16-24 (9 rows). The new data set only
# has 7 rows, so remove two rows from the middle of the table. This keeps
# the specially-bordered closing row (originally row 24) intact as the new
# last row of the table (row 22) once the rows above it shift up.
$ws.Rows("20:21").Delete()

# Header summary fields
$ws.Range("E11").Value = 305533
$ws.Range("C13").Value = 3
$ws.Range("F13").Value = 6

# Row 16: KAROL ANDREA PEÑATE GARCIA
$ws.Range("C16").Value = "1143382362"
$ws.Range("D16").Value = "KAROL ANDREA PEÑATE GARCIA"
$ws.Range("E16").Value = "1812"
$ws.Range("F16").Value = 14583
$ws.Range("G16").Value = 877803

# Row 17: HAROLD MANUEL BARRIOS PAJARO
$ws.Range("C17").Value = "1143362003"
$ws.Range("D17").Value = "HAROLD MANUEL BARRIOS PAJARO"
$ws.Range("E17").Value = "1812"
$ws.Range("F17").Value = 6250
$ws.Range("G17").Value = 828116

# Rows 18-22: YAMIRIS LARA MENDOZA, one row per period 2504-2508
$ws.Range("C18").Value = "1047422669"
$ws.Range("D18").Value = "YAMIRIS LARA MENDOZA"
$ws.Range("E18").Value = "2504"
$ws.Range("F18").Value = 56940
$ws.Range("G18").Value = 1423500

$ws.Range("C19").Value = "1047422669"
$ws.Range("D19").Value = "YAMIRIS LARA MENDOZA"
$ws.Range("E19").Value = "2505"
$ws.Range("F19").Value = 56940
$ws.Range("G19").Value = 1423500

$ws.Range("C20").Value = "1047422669"
$ws.Range("D20").Value = "YAMIRIS LARA MENDOZA"
$ws.Range("E20").Value = "2506"
$ws.Range("F20").Value = 56940
$ws.Range("G20").Value = 1423500

$ws.Range("C21").Value = "1047422669"
$ws.Range("D21").Value = "YAMIRIS LARA MENDOZA"
$ws.Range("E21").Value = "2507"
$ws.Range("F21").Value = 56940
$ws.Range("G21").Value = 1423500

$ws.Range("C22").Value = "1047422669"
$ws.Range("D22").Value = "YAMIRIS LARA MENDOZA"
$ws.Range("E22").Value = "2508"
$ws.Range("F22").Value = 56940
$ws.Range("G22").Value = 1423500
